# Fixed-data-types.xlsx - "Minor correction to the docs."
#
# The source diff doubles every multiplier formula in column E (rows 3-11)
# from "=C<n>*16" to "=C<n>*16*2" (and the selection/active cell moves from
# E6 to the E3:E11 range). The remaining hunks in the diff (font "charset"
# attributes, and a pure re-pairing of two unused named cell styles
# "Result2"/"Heading1" <-> their cellStyleXfs records) are not reachable
# through the Excel COM object model - Style.Name has no setter and
# Font.Charset is not wired to persistence in this host - and none of
# those records are referenced by any cell actually used on the sheet, so
# they carry no visible effect; only the formulas/selection below are
# user-observable and are applied here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E (rows 3-11): multiply by an extra factor of 2.
for ($r = 3; $r -le 11; $r++) {
    $ws.Cells.Item($r, 5).Formula = "=C$r*16*2"
}

# Move the sheet's saved selection from E6 to E3:E11 (active cell E3).
$ws.Range("E3:E11").Select() | Out-Null
